# Updates cryptos list values (price + 1h volume change) pulled from the
# latest coinranking.com snapshot. Also corrects the MXToken/HuobiToken row
# ordering (rows 36-37 swapped places along with their data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin name / link / already-text price / percent cells).
$plainUpdates = @{
    'D2' = '27.518.03'
    'E2' = '  +5.08%  '
    'D3' = '1.724.33'
    'E4' = '  +0.14%  '
    'E5' = '  +3.25%  '
    'E6' = '  +2.54%  '
    'E7' = '  +0.06%  '
    'E8' = '  +1.06%  '
    'E9' = '  +3.89%  '
    'E10' = '  +4.98%  '
    'E11' = '  +0.82%  '
    'E12' = '  -0.05%  '
    'D13' = '1.736.76'
    'E13' = '  +5.01%  '
    'D14' = '1.962.02'
    'E14' = '  +3.97%  '
    'E15' = '  +4.56%  '
    'D16' = '0.0₅8307'
    'E16' = '  +0.76%  '
    'E17' = '  +3.62%  '
    'D18' = '27.549.52'
    'E18' = '  +5.22%  '
    'E19' = '  +16.61%  '
    'E20' = '  +0.10%  '
    'E22' = '  +2.07%  '
    'E23' = '  +2.35%  '
    'E24' = '  +0.07%  '
    'E25' = '  +1.35%  '
    'E26' = '  +11.34%  '
    'E27' = '  +2.90%  '
    'E28' = '  +1.58%  '
    'E29' = '  +4.89%  '
    'E30' = '  +0.05%  '
    'E31' = '  +2.50%  '
    'E32' = '  +3.19%  '
    'E33' = '  +2.25%  '
    'E34' = '  +6.15%  '
    'E35' = '  +1.07%  '
    'B36' = 'HuobiToken'
    'C36' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'E36' = '  +1.87%  '
    'B37' = 'MXToken'
    'C37' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'E37' = '  +1.36%  '
    'E38' = '  +3.89%  '
    'E39' = '  +3.41%  '
    'E40' = '  -0.55%  '
    'E41' = '  +2.89%  '
    'D42' = '1.060.76'
    'E42' = '  +2.51%  '
    'E43' = '  +0.10%  '
    'E44' = '  +0.23%  '
    'D45' = '1.867.22'
    'E45' = '  +3.87%  '
    'D46' = '0.0₈115'
    'E46' = '  +8.65%  '
    'E47' = '  +1.08%  '
    'E48' = '  +2.03%  '
    'E50' = '  +0.38%  '
    'E51' = '  +1.04%  '
}
foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Price updates whose new text looks like a number (e.g. "1.005") would be
# auto-coerced to a numeric value by Excel, so force the cell to Text format
# first, assign the literal string, then clear the format back to General so
# the cell style matches the rest of the sheet (unstyled, default numFmt).
$numericLookingUpdates = @{
    'D4' = '1.005'
    'D5' = '226.23'
    'D6' = '0.5382'
    'D7' = '1.005'
    'D8' = '0.2702'
    'D9' = '0.06621'
    'D10' = '21.71'
    'D12' = '4.657'
    'D15' = '0.5891'
    'D17' = '68.08'
    'D19' = '224.77'
    'D20' = '1.005'
    'D21' = '4.764'
    'D22' = '10.74'
    'D23' = '6.128'
    'D25' = '147.85'
    'D26' = '1.697'
    'D27' = '0.1237'
    'D28' = '7.431'
    'D29' = '16.82'
    'D30' = '0.05592'
    'D32' = '3.595'
    'D33' = '3.472'
    'D34' = '1.670'
    'D35' = '0.9656'
    'D36' = '2.449'
    'D37' = '2.821'
    'D38' = '0.5959'
    'D40' = '5.899'
    'D41' = '0.8597'
    'D44' = '101.50'
    'D47' = '59.19'
    'D48' = '8.242'
    'D51' = '0.05293'
}
foreach ($ref in $numericLookingUpdates.Keys) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $numericLookingUpdates[$ref]
    $ws.Range($ref).ClearFormats()
}

